$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 6490.7
$ws.Range("J28").Value = 14724.75
$ws.Range("L28").Value = 14724.75
$ws.Range("N28").Value = -15694.75

# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 225.38889
$ws.Range("I33").Value = 237.9375
$ws.Range("K33").Value = 237.9375
$ws.Range("M33").Value = -8.9375

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 125035150
$ws.Range("J40").Value = 166713040
$ws.Range("L40").Value = 166713040
$ws.Range("N40").Value = -166713390

# Row 76 (Leve Item ID 12602)
$ws.Range("H76").Value = 6306.0713
$ws.Range("I76").Value = 8463.875
$ws.Range("K76").Value = 8463.875
$ws.Range("M76").Value = -8148.875

# Row 79 (Leve Item ID 12602)
$ws.Range("H79").Value = 6306.0713
$ws.Range("I79").Value = 8463.875
$ws.Range("K79").Value = 8463.875
$ws.Range("M79").Value = -7371.875

# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 1076.6923
$ws.Range("J107").Value = 2824.1428
$ws.Range("L107").Value = 2824.1428
$ws.Range("N107").Value = -6664.1428

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 16042.267
$ws.Range("I116").Value = 9998.5
$ws.Range("J116").Value = 16972.076
$ws.Range("K116").Value = 9998.5
$ws.Range("L116").Value = 16972.076
$ws.Range("M116").Value = -6556.5
$ws.Range("N116").Value = -23856.076

# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 2100
$ws.Range("I127").Value = 2100
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 6300
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -1340
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 30 (Leve Item ID 2712)
$ws.Range("H30").Value = 750005000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 750005000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 750005000
$ws.Range("N30").Value = -750005300
$ws.Range("M30").ClearContents()

# Row 103 (Leve Item ID 18533)
$ws.Range("H103").Value = 92892
$ws.Range("J103").Value = 92892
$ws.Range("L103").Value = 92892
$ws.Range("N103").Value = -95236

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 2262.6191
$ws.Range("I20").Value = 1894.1034
$ws.Range("J20").Value = 3084.6924
$ws.Range("K20").Value = 1894.1034
$ws.Range("L20").Value = 3084.6924
$ws.Range("M20").Value = -1647.1034
$ws.Range("N20").Value = -3578.6924

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 66.1579
$ws.Range("J7").Value = 58.77778
$ws.Range("L7").Value = 58.77778
$ws.Range("N7").Value = -284.77778

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 2727.2778
$ws.Range("I58").Value = 2006.2667
$ws.Range("J58").Value = 6332.3335
$ws.Range("K58").Value = 2006.2667
$ws.Range("L58").Value = 6332.3335
$ws.Range("M58").Value = -1803.2667
$ws.Range("N58").Value = -6738.3335

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1951.4359
$ws.Range("I134").Value = 1632.7333
$ws.Range("K134").Value = 4898.199900000001
$ws.Range("M134").Value = -2363.199900000001

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 2727.2778
$ws.Range("I136").Value = 2006.2667
$ws.Range("J136").Value = 6332.3335
$ws.Range("K136").Value = 6018.800099999999
$ws.Range("L136").Value = 18997.0005
$ws.Range("M136").Value = -3468.800099999999
$ws.Range("N136").Value = -24097.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 274
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 288.92307
$ws.Range("K2").Value = 480
$ws.Range("L2").Value = 1733.53842
$ws.Range("M2").Value = -367
$ws.Range("N2").Value = -1959.53842

# Row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 166674660
$ws.Range("I11").Value = 250000900
$ws.Range("K11").Value = 750002700
$ws.Range("M11").Value = -750002560

# Row 26 (Leve Item ID 4746)
$ws.Range("H26").Value = 5630.143
$ws.Range("I26").Value = 1075.8
$ws.Range("J26").Value = 17016
$ws.Range("K26").Value = 3227.4
$ws.Range("L26").Value = 51048
$ws.Range("M26").Value = -2939.4
$ws.Range("N26").Value = -51624

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 14266.182
$ws.Range("I70").Value = 12864
$ws.Range("K70").Value = 12864
$ws.Range("M70").Value = -12594

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 14266.182
$ws.Range("I73").Value = 12864
$ws.Range("K73").Value = 12864
$ws.Range("M73").Value = -11928

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3239062.8
$ws.Range("I132").Value = 3526.7827
$ws.Range("J132").Value = 14687882
$ws.Range("K132").Value = 10580.3481
$ws.Range("L132").Value = 44063646
$ws.Range("M132").Value = -8050.348100000001
$ws.Range("N132").Value = -44068706

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 6739.838
$ws.Range("I7").Value = 6457.1665
$ws.Range("K7").Value = 6457.1665
$ws.Range("M7").Value = -6345.1665

# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 2131.32
$ws.Range("I16").Value = 2089
$ws.Range("J16").Value = 2300.6
$ws.Range("K16").Value = 2089
$ws.Range("L16").Value = 2300.6
$ws.Range("M16").Value = -1919
$ws.Range("N16").Value = -2640.6

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 11498.23
$ws.Range("I22").Value = 38996.668
$ws.Range("K22").Value = 38996.668
$ws.Range("M22").Value = -38701.668

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 11498.23
$ws.Range("I27").Value = 38996.668
$ws.Range("K27").Value = 38996.668
$ws.Range("M27").Value = -38889.668

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 4251.08
$ws.Range("I40").Value = 3285.4736
$ws.Range("J40").Value = 7308.8335
$ws.Range("K40").Value = 3285.4736
$ws.Range("L40").Value = 7308.8335
$ws.Range("M40").Value = -3149.4736
$ws.Range("N40").Value = -7580.8335

# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 15643706
$ws.Range("I100").Value = 2445.25
$ws.Range("K100").Value = 2445.25
$ws.Range("M100").Value = -1904.25

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 6739.838
$ws.Range("I126").Value = 6457.1665
$ws.Range("K126").Value = 19371.4995
$ws.Range("M126").Value = -16901.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 69 (Leve Item ID 10951)
$ws.Range("H69").Value = 14999
$ws.Range("J69").Value = 14999
$ws.Range("L69").Value = 14999
$ws.Range("N69").Value = -16497

# Row 72 (Leve Item ID 10951)
$ws.Range("H72").Value = 14999
$ws.Range("J72").Value = 14999
$ws.Range("L72").Value = 44997
$ws.Range("N72").Value = -52485

# Row 105 (Leve Item ID 18710)
$ws.Range("H105").Value = 34997
$ws.Range("J105").Value = 34997
$ws.Range("L105").Value = 34997
$ws.Range("N105").Value = -41985

# Row 140 (Leve Item ID 42506)
$ws.Range("H140").Value = 53729.168
$ws.Range("J140").Value = 53729.168
$ws.Range("L140").Value = 53729.168
$ws.Range("N140").Value = -64089.168
